# Update the "Förändrad" (Changed) date column (C) for rows 2-29
# from serial date 45576 (2024-10-11) to 45577 (2024-10-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45576) {
        $cell.Value2 = 45577
    }
}
